$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to remain text (avoid Excel auto-numeric conversion)
$dRange = $ws.Range('D2:D51')
$dRange.NumberFormat = '@'

$ws.Range('D2').Value = '28.300.94'
$ws.Range('E2').Value = '  -0.01%  '
$ws.Range('D3').Value = '1.864.56'
$ws.Range('E3').Value = '  +2.95%  '
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '310.77'
$ws.Range('E5').Value = '  -0.53%  '
$ws.Range('D6').Value = '0.9998'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.4985'
$ws.Range('E7').Value = '  -3.42%  '
$ws.Range('D8').Value = '0.3975'
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = '0.09988'
$ws.Range('E9').Value = '  +27.60%  '
$ws.Range('D10').Value = '1.121'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('D11').Value = '41.40'
$ws.Range('E11').Value = '  +0.88%  '
$ws.Range('D12').Value = '6.487'
$ws.Range('E12').Value = '  +2.14%  '
$ws.Range('D13').Value = '20.93'
$ws.Range('E13').Value = '  +2.44%  '
$ws.Range('D14').Value = '1.866.94'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '1.0000'
$ws.Range('E15').Value = '  -0.18%  '
$ws.Range('D16').Value = '7.399'
$ws.Range('D17').Value = '0.00001146'
$ws.Range('E17').Value = '  +5.66%  '
$ws.Range('D18').Value = '93.55'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D19').Value = '0.06644'
$ws.Range('E19').Value = '  +0.96%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '17.38'
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('B21').Value = 'Dai'
$ws.Range('C21').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D21').Value = '1.000'
$ws.Range('E21').Value = '  -0.14%  '
$ws.Range('D22').Value = '6.072'
$ws.Range('E22').Value = '  +1.10%  '
$ws.Range('D23').Value = '28.377.32'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').Value = '11.35'
$ws.Range('E24').Value = '  +2.05%  '
$ws.Range('D25').Value = '2.246'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('D26').Value = '21.21'
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('D27').Value = '2.079.00'
$ws.Range('E27').Value = '  +2.96%  '
$ws.Range('D28').Value = '2.499'
$ws.Range('E28').Value = '  +3.09%  '
$ws.Range('D29').Value = '157.41'
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('D30').Value = '127.81'
$ws.Range('E30').Value = '  +0.26%  '
$ws.Range('D31').Value = '0.1058'
$ws.Range('E31').Value = '  -3.78%  '
$ws.Range('D32').Value = '1.052'
$ws.Range('E32').Value = '  -0.91%  '
$ws.Range('D33').Value = '5.634'
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').Value = '3.598'
$ws.Range('E34').Value = '  -1.66%  '
$ws.Range('D35').Value = '0.06809'
$ws.Range('D36').Value = '9.237'
$ws.Range('E36').Value = '  +1.24%  '
$ws.Range('D37').Value = '0.02383'
$ws.Range('E37').Value = '  +0.99%  '
$ws.Range('D38').Value = '0.2173'
$ws.Range('E38').Value = '  -0.72%  '
$ws.Range('D39').Value = '5.024'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('D40').Value = '11.50'
$ws.Range('E40').Value = '  -0.55%  '
$ws.Range('D41').Value = '0.6301'
$ws.Range('E41').Value = '  +1.94%  '
$ws.Range('D42').Value = '1.180'
$ws.Range('E42').Value = '  +1.95%  '
$ws.Range('D43').Value = '0.9995'
$ws.Range('E43').Value = '  -0.16%  '
$ws.Range('D44').Value = '13.44'
$ws.Range('E44').Value = '  +1.13%  '
$ws.Range('D45').Value = '0.6000'
$ws.Range('E45').Value = '  -0.04%  '
$ws.Range('D46').Value = '1.283'
$ws.Range('E46').Value = '  -1.74%  '
$ws.Range('D47').Value = '3.668'
$ws.Range('E47').Value = '  -1.98%  '
$ws.Range('D48').Value = '125.05'
$ws.Range('E48').Value = '  -0.12%  '
$ws.Range('D49').Value = '1.993'
$ws.Range('E49').Value = '  +3.59%  '
$ws.Range('D50').Value = '1.189'
$ws.Range('E50').Value = '  -2.37%  '
$ws.Range('D51').Value = '1.119'
$ws.Range('E51').Value = '  +4.07%  '

# Restore default style on column D so no stray style index is introduced
$dRange.Style = 'Normal'

